$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '25.836.83'
Set-TextValue $ws.Range('E2') '  -0.20%  '
Set-TextValue $ws.Range('D3') '1.639.05'
Set-TextValue $ws.Range('E3') '  -0.04%  '
Set-TextValue $ws.Range('E4') '  -0.30%  '
Set-TextValue $ws.Range('D5') '215.58'
Set-TextValue $ws.Range('E5') '  -0.04%  '
Set-TextValue $ws.Range('E6') '  -0.52%  '
Set-TextValue $ws.Range('E7') '  -0.29%  '
Set-TextValue $ws.Range('E8') '  -0.79%  '
Set-TextValue $ws.Range('E9') '  -1.07%  '
Set-TextValue $ws.Range('D10') '19.83'
Set-TextValue $ws.Range('E10') '  -2.02%  '
Set-TextValue $ws.Range('E11') '  +1.59%  '
Set-TextValue $ws.Range('D13') '1.864.37'
Set-TextValue $ws.Range('E13') '  -0.08%  '
Set-TextValue $ws.Range('D14') '1.639.38'
Set-TextValue $ws.Range('E14') '  -0.40%  '
Set-TextValue $ws.Range('D15') '0.564'
Set-TextValue $ws.Range('E15') '  -0.41%  '
Set-TextValue $ws.Range('E16') '  -0.11%  '
Set-TextValue $ws.Range('D17') '63.31'
Set-TextValue $ws.Range('E17') '  -0.26%  '
Set-TextValue $ws.Range('D18') '25.851.24'
Set-TextValue $ws.Range('E18') '  -0.22%  '
Set-TextValue $ws.Range('E20') '  +2.02%  '
Set-TextValue $ws.Range('D21') '193.12'
Set-TextValue $ws.Range('E21') '  -0.78%  '
Set-TextValue $ws.Range('D22') '10.02'
Set-TextValue $ws.Range('E22') '  +0.30%  '
Set-TextValue $ws.Range('D23') '6.36'
Set-TextValue $ws.Range('E23') '  +1.82%  '
Set-TextValue $ws.Range('B24') 'Toncoin'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D24') '1.82'
Set-TextValue $ws.Range('E24') '  +3.80%  '
Set-TextValue $ws.Range('B25') 'BinanceUSD'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D25') '1.00'
Set-TextValue $ws.Range('E25') '  -0.26%  '
Set-TextValue $ws.Range('D26') '142.81'
Set-TextValue $ws.Range('E26') '  +3.06%  '
Set-TextValue $ws.Range('D28') '6.98'
Set-TextValue $ws.Range('E28') '  +1.72%  '
Set-TextValue $ws.Range('E29') '  +0.08%  '
Set-TextValue $ws.Range('E30') '  -0.53%  '
Set-TextValue $ws.Range('E31') '  -0.24%  '
Set-TextValue $ws.Range('E32') '  +1.10%  '
Set-TextValue $ws.Range('E33') '  -0.49%  '
Set-TextValue $ws.Range('D34') '1.59'
Set-TextValue $ws.Range('E34') '  +0.63%  '
Set-TextValue $ws.Range('E35') '  -0.42%  '
Set-TextValue $ws.Range('D36') '0.910'
Set-TextValue $ws.Range('E36') '  +0.01%  '
Set-TextValue $ws.Range('D37') '1.133.46'
Set-TextValue $ws.Range('E37') '  +0.53%  '
Set-TextValue $ws.Range('E38') '  -1.50%  '
Set-TextValue $ws.Range('D39') '0.547'
Set-TextValue $ws.Range('E39') '  -1.21%  '
Set-TextValue $ws.Range('E40') '  -0.43%  '
Set-TextValue $ws.Range('E41') '  +0.07%  '
Set-TextValue $ws.Range('D42') '5.56'
Set-TextValue $ws.Range('E42') '  +0.90%  '
Set-TextValue $ws.Range('D43') '100.41'
Set-TextValue $ws.Range('E43') '  +0.82%  '
Set-TextValue $ws.Range('E44') '  +0.66%  '
Set-TextValue $ws.Range('D45') '1.773.52'
Set-TextValue $ws.Range('E45') '  -0.27%  '
Set-TextValue $ws.Range('E46') '  +3.56%  '
Set-TextValue $ws.Range('D47') '55.44'
Set-TextValue $ws.Range('E47') '  -0.42%  '
Set-TextValue $ws.Range('E48') '  -1.39%  '
Set-TextValue $ws.Range('E49') '  -0.16%  '
Set-TextValue $ws.Range('D50') '1.43'
Set-TextValue $ws.Range('E50') '  +4.23%  '
Set-TextValue $ws.Range('E51') '  +3.50%  '
